# Auto-generated edit script: updates leve profit figures across
# multiple Sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) per the
# scheduled market-price refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2160.0557
$ws.Range("I33").Value = 1498.3846
$ws.Range("K33").Value = 1498.3846
$ws.Range("M33").Value = -1269.3846
$ws.Range("H98").Value = 31828108
$ws.Range("I98").Value = 15386815
$ws.Range("J98").Value = 55576644
$ws.Range("K98").Value = 15386815
$ws.Range("L98").Value = 55576644
$ws.Range("M98").Value = -15385317
$ws.Range("N98").Value = -55579640
$ws.Range("H122").Value = 31828108
$ws.Range("I122").Value = 15386815
$ws.Range("J122").Value = 55576644
$ws.Range("K122").Value = 46160445
$ws.Range("L122").Value = 166729932
$ws.Range("M122").Value = -46157995
$ws.Range("N122").Value = -166734832
$ws.Range("H138").Value = 2979
$ws.Range("I138").Value = 2962.2856
$ws.Range("J138").Value = 2990.3225
$ws.Range("K138").Value = 8886.856800000001
$ws.Range("L138").Value = 8970.967500000001
$ws.Range("M138").Value = -3746.856800000001
$ws.Range("N138").Value = -19250.9675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1552.9333
$ws.Range("I122").Value = 710
$ws.Range("J122").Value = 1974.4
$ws.Range("K122").Value = 2130
$ws.Range("L122").Value = 5923.200000000001
$ws.Range("M122").Value = 320
$ws.Range("N122").Value = -10823.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4778862.5
$ws.Range("I58").Value = 2385541.2
$ws.Range("J58").Value = 22728772
$ws.Range("K58").Value = 2385541.2
$ws.Range("L58").Value = 22728772
$ws.Range("M58").Value = -2385338.2
$ws.Range("N58").Value = -22729178
$ws.Range("H88").Value = 26468.6
$ws.Range("J88").Value = 26468.6
$ws.Range("L88").Value = 26468.6
$ws.Range("N88").Value = -27280.6
$ws.Range("H91").Value = 26468.6
$ws.Range("J91").Value = 26468.6
$ws.Range("L91").Value = 26468.6
$ws.Range("N91").Value = -29276.6
$ws.Range("H134").Value = 1673478.6
$ws.Range("I134").Value = 7295.4707
$ws.Range("J134").Value = 5719923.5
$ws.Range("K134").Value = 21886.4121
$ws.Range("L134").Value = 17159770.5
$ws.Range("M134").Value = -19351.4121
$ws.Range("N134").Value = -17164840.5
$ws.Range("H136").Value = 4778862.5
$ws.Range("I136").Value = 2385541.2
$ws.Range("J136").Value = 22728772
$ws.Range("K136").Value = 7156623.600000001
$ws.Range("L136").Value = 68186316
$ws.Range("M136").Value = -7154073.600000001
$ws.Range("N136").Value = -68191416

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3676090.8
$ws.Range("I5").Value = 4274221
$ws.Range("J5").Value = 2778895.5
$ws.Range("K5").Value = 12822663
$ws.Range("L5").Value = 8336686.5
$ws.Range("M5").Value = -12822551
$ws.Range("N5").Value = -8336910.5
$ws.Range("H135").Value = 3676090.8
$ws.Range("I135").Value = 4274221
$ws.Range("J135").Value = 2778895.5
$ws.Range("K135").Value = 38467989
$ws.Range("L135").Value = 25010059.5
$ws.Range("M135").Value = -38465454
$ws.Range("N135").Value = -25015129.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2769596
$ws.Range("I70").Value = 1689229
$ws.Range("J70").Value = 4552201.5
$ws.Range("K70").Value = 1689229
$ws.Range("L70").Value = 4552201.5
$ws.Range("M70").Value = -1688959
$ws.Range("N70").Value = -4552741.5
$ws.Range("H73").Value = 2769596
$ws.Range("I73").Value = 1689229
$ws.Range("J73").Value = 4552201.5
$ws.Range("K73").Value = 1689229
$ws.Range("L73").Value = 4552201.5
$ws.Range("M73").Value = -1688293
$ws.Range("N73").Value = -4554073.5
$ws.Range("H102").Value = 6213.8823
$ws.Range("I102").Value = 7095.4287
$ws.Range("J102").Value = 2100
$ws.Range("K102").Value = 7095.4287
$ws.Range("L102").Value = 2100
$ws.Range("M102").Value = -5473.4287
$ws.Range("N102").Value = -5344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1715.1538
$ws.Range("I7").Value = 1071.4286
$ws.Range("J7").Value = 2466.1667
$ws.Range("K7").Value = 1071.4286
$ws.Range("L7").Value = 2466.1667
$ws.Range("M7").Value = -959.4286
$ws.Range("N7").Value = -2690.1667
$ws.Range("H40").Value = 1698.25
$ws.Range("I40").Value = 976.61536
$ws.Range("J40").Value = 3038.4285
$ws.Range("K40").Value = 976.61536
$ws.Range("L40").Value = 3038.4285
$ws.Range("M40").Value = -840.61536
$ws.Range("N40").Value = -3310.4285
$ws.Range("H62").Value = 35111
$ws.Range("J62").Value = 35111
$ws.Range("L62").Value = 35111
$ws.Range("N62").Value = -36359
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 28330
$ws.Range("J64").Value = 24995
$ws.Range("L64").Value = 24995
$ws.Range("N64").Value = -25445
$ws.Range("H65").Value = 35111
$ws.Range("J65").Value = 35111
$ws.Range("L65").Value = 105333
$ws.Range("N65").Value = -111573
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 28330
$ws.Range("J67").Value = 24995
$ws.Range("L67").Value = 24995
$ws.Range("N67").Value = -26555
$ws.Range("H68").Value = 2465.5715
$ws.Range("J68").Value = 2920
$ws.Range("L68").Value = 2920
$ws.Range("N68").Value = -4418
$ws.Range("H71").Value = 2465.5715
$ws.Range("J71").Value = 2920
$ws.Range("L71").Value = 14600
$ws.Range("N71").Value = -22088
$ws.Range("H122").Value = 9222408
$ws.Range("I122").Value = 1184599.1
$ws.Range("J122").Value = 33335834
$ws.Range("K122").Value = 3553797.3
$ws.Range("L122").Value = 100007502
$ws.Range("M122").Value = -3551347.3
$ws.Range("N122").Value = -100012402
$ws.Range("H126").Value = 1715.1538
$ws.Range("I126").Value = 1071.4286
$ws.Range("J126").Value = 2466.1667
$ws.Range("K126").Value = 3214.2858
$ws.Range("L126").Value = 7398.500100000001
$ws.Range("M126").Value = -744.2857999999997
$ws.Range("N126").Value = -12338.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982
$ws.Range("H122").Value = 1348.7407
$ws.Range("I122").Value = 1302.1333
$ws.Range("J122").Value = 1407
$ws.Range("K122").Value = 3906.3999
$ws.Range("L122").Value = 4221
$ws.Range("M122").Value = -1456.3999
$ws.Range("N122").Value = -9121
$ws.Range("H126").Value = 31252614
$ws.Range("I126").Value = 83333840
$ws.Range("J126").Value = 3879.8
$ws.Range("K126").Value = 250001520
$ws.Range("L126").Value = 11639.4
$ws.Range("M126").Value = -249999050
$ws.Range("N126").Value = -16579.4
